$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new worksheet column at N (14) - shifts "link" (N->O) and "data_path" (O->P)
$ws.Columns.Item(14).Insert()

# Header for the new column
$ws.Range("N1").Value = "delivery_format"

# Re-stamp headers for the columns that shifted right, so the table can
# pick up their names again at their new position
$ws.Range("O1").Value = "link"
$ws.Range("P1").Value = "data_path"

# Data for the new delivery_format column (xlsx must be written before csv so
# the shared-string table gets the same ordering as the target workbook)
$ws.Range("N3").Value = "xlsx"
$ws.Range("N2").Value = "csv"
$ws.Range("N4").Value = "csv"

# Grow the table to include the new column
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:P4"))

# Column width for the new column (matches the narrower delivery_format header)
$ws.Columns.Item(14).ColumnWidth = 15.166666666666666

$sel = $ws.Range("N7")
$sel.Select()
